{"js": "// The document has two single-run paragraphs whose text needs to be split\n// into several <w:t> segments joined by manual line breaks (<w:br/>),\n// without altering the run formatting (both halves stay in the same run).\n//\n// 1) In the \"Avalia\u00e7\u00e3o\" paragraph, the \"Crit\u00e9rio:\" sentence gets a line\n//    break right before \"NF = (P1 + P2 + T)/3...\".\n// 2) In the \"Bibliografia\" paragraph, a line break is inserted before each\n//    of the four later numbered references (\"2.\", \"3.\", \"4.\", \"5.\"), so the\n//    five-item list ends up one reference per line.\n//\n// Word represents a manual line break as a run child \"\\u000b\" (vertical\n// tab) when read/written through the text APIs, so inserting that\n// character right before each target phrase reproduces the <w:t>/<w:br/>\n// split seen in the diff.\n\nconst body = context.document.body;\n\nasync function insertLineBreakBefore(searchText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for ${JSON.stringify(searchText)}, found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(\"\\u000b\", Word.InsertLocation.before);\n  await context.sync();\n}\n\n// --- Change 1: \"Crit\u00e9rio\" paragraph \u2014 break before the \"NF = ...\" formula.\nawait insertLineBreakBefore(\"NF = (P1 + P2 + T)/3\");\n\n// --- Change 2: \"Bibliografia\" paragraph \u2014 break before each later entry.\nawait insertLineBreakBefore(\"2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL\");\nawait insertLineBreakBefore(\"3. HIBBELER, R.C.\");\nawait insertLineBreakBefore(\"4. MERIAM, J.L. KRAIGE\");\nawait insertLineBreakBefore(\"5. RUIZ, C.C.de La P.\");\n", "ps1": "# The document has two single-run paragraphs whose text needs to be split\n# into several runs joined by manual line breaks, without altering the run\n# formatting (the text before and after each break stays in one run, just\n# like the original XML shows <w:t>...</w:t><w:br/><w:t>...</w:t>).\n#\n# 1) In the \"Avaliacao\" paragraph, the \"Criterio:\" sentence gets a line\n#    break right before \"NF = (P1 + P2 + T)/3...\".\n# 2) In the \"Bibliografia\" paragraph, a line break is inserted before each\n#    of the four later numbered references (\"2.\", \"3.\", \"4.\", \"5.\"), so the\n#    five-item list ends up one reference per line.\n#\n# A manual line break in Word is the chr(11) vertical-tab character; the\n# helper below finds the target phrase, collapses the found range to its\n# start, and inserts that character immediately before it.\n\n$d = $word.ActiveDocument\n\nfunction Insert-LineBreak-Before([string]$needle) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $needle\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $needle\"\n    }\n    $rng.Collapse(1)  # wdCollapseStart\n    $rng.InsertBefore([char]11)\n}\n\n# --- Change 1: \"Criterio\" paragraph - break before the \"NF = ...\" formula.\nInsert-LineBreak-Before \"NF = (P1 + P2 + T)/3\"\n\n# --- Change 2: \"Bibliografia\" paragraph - break before each later entry.\nInsert-LineBreak-Before \"2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL\"\nInsert-LineBreak-Before \"3. HIBBELER, R.C.\"\nInsert-LineBreak-Before \"4. MERIAM, J.L. KRAIGE\"\nInsert-LineBreak-Before \"5. RUIZ, C.C.de La P.\"\n"}
